$d = $word.ActiveDocument
$find = $d.Content.Find
$ok = $find.Execute("베타 관계", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "Could not find anchor text '베타 관계'"
}
$anchor = $find.Parent.Duplicate
$insertionPoint = $d.Range($anchor.End, $anchor.End)
$xml = '<w:p /><w:p><w:pPr><w:rPr><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t xml:space="preserve">033 </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t>E-R</w:t></w:r><w:r><w:rPr><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t>개체-관계)</w:t></w:r><w:r><w:rPr><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t>모델</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">- 개체와 개체 간이 관계를 기본 요소로 이용하여 현실 세계의 </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">무질서한 데이터를 개념적인 논리 </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:lastRenderedPageBreak /><w:t>데이터로 표현하기 위한 방법</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>데이터를 개체,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>관계,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>속성으로 묘사</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>E-R</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>다이어그램 :</w:t></w:r><w:proofErr w:type="gramEnd" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve"> 사각형 </w:t></w:r><w:r><w:t xml:space="preserve">/ </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">마름모 </w:t></w:r><w:r><w:t xml:space="preserve">/ </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">타원 </w:t></w:r><w:r><w:t xml:space="preserve">/ </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">이중 타원 </w:t></w:r><w:r><w:t xml:space="preserve">/ </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">밑줄 타원 </w:t></w:r><w:r><w:t xml:space="preserve">/ </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">복수 타원 </w:t></w:r><w:r><w:t xml:space="preserve">/ </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">관계 </w:t></w:r><w:r><w:t xml:space="preserve">/ </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>선,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>링크</w:t></w:r></w:p><w:p /><w:p><w:pPr><w:rPr><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia" /><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t xml:space="preserve">034 관계형 데이터베이스의 구조 </w:t></w:r><w:r><w:rPr><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t xml:space="preserve">/ </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /><w:sz w:val="22" /><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF" /></w:rPr><w:t>관계형 데이터 모델</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">관계형 데이터베이스 </w:t></w:r><w:r><w:t>- 2</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>차원적인 표를 이용해서 데이터 상호 관계를 정의하는 데이터베이스</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart" /><w:proofErr w:type="gramStart" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>튜플</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>:</w:t></w:r><w:proofErr w:type="gramEnd" /><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>릴레이션을</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve"> 구성하는 각각의 행</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>속성 :</w:t></w:r><w:proofErr w:type="gramEnd" /><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>데이터베이스를 구성하는 가장 작은 논리적 단위</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">도메인 </w:t></w:r><w:r><w:t>:</w:t></w:r><w:proofErr w:type="gramEnd" /><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">하나의 </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>애트리뷰트가</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve"> 취할 수 있는 같은 타입의 </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>원자값들의</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve"> 집합</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">관계형 데이터 모델 </w:t></w:r><w:r><w:t>- 2</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t xml:space="preserve">차원적인 표를 이용해서 데이터 상호 관계를 정의하는 </w:t></w:r><w:r><w:t>DB</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" /></w:rPr><w:t>구조를 말함</w:t></w:r></w:p>'
$insertionPoint.InsertXML($xml)
